$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("F1").Value = 0.03
$ws.Range("F2").Value = 0.01
$ws.Range("I2").Value = -0.03
$ws.Range("I6").Value = 0.06

# Update the view: scroll back to A1 (remove topLeftCell="B1") and change selection to I3
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I3").Select()
